# Insert a new row at row 270, shifting the existing rows 270-372 down to 271-373,
# then populate the new row 270 with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 270 (pushes old 270..372 down to 271..373)
$ws.Rows.Item(270).Insert()

# Fill in the new row 270 with the new data record.
$ws.Range("A270").Value = 11
$ws.Range("B270").Value = "Vega Monumental Concepción"
$ws.Range("C270").Value = "Bíobío"
$ws.Range("D270").Value = 44924
$ws.Range("E270").Value = 8
$ws.Range("F270").Value = 100114001
$ws.Range("G270").Value = "Papa"
$ws.Range("H270").Value = "Patagonia"
$ws.Range("I270").Value = "1a (cosecha)"
$ws.Range("J270").Value = 5000
$ws.Range("K270").Value = 11000
$ws.Range("L270").Value = 12000
$ws.Range("M270").Value = 11500
$ws.Range("N270").Value = "`$/saco 25 kilos"
$ws.Range("O270").Value = "Provincia de Arauco"
$ws.Range("P270").Value = 460
$ws.Range("Q270").Value = 25
$ws.Range("R270").Value = "Hortaliza"
